$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "267.47"
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "6.329"
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "0.06194"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "3.592"
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "6.674"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "1.386"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.8287"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.01361"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.08235"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.03396"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.03160"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.09301"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "3.915"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.001717"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.04857"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.006247"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.005386"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "0.001093"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.0001503"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "3.761"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.370"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.3346"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.1216"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.0002688"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.04668"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.006901"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.1153"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.003469"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.01227"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.00006275"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.00000000752"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.7906"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.1579"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.00002105"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.01243"
